$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fase 2")
$ws.Range("A28").Value = "test"
